# Insert a new data row right before the current row 56 (shifts rows
# 56..129 down to 57..130, dimension grows from A1:R129 to A1:R130),
# then populate the newly inserted row 56 with its own values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(56).Insert()

$ws.Range("A56").Value = 3
$ws.Range("B56").Value = "Femacal de La Calera"
$ws.Range("C56").Value = "Coquimbo"
$ws.Range("D56").Value = 44483
$ws.Range("E56").Value = 5
$ws.Range("F56").Value = 100112010
$ws.Range("G56").Value = "Achicoria"
$ws.Range("H56").Value = "Sin especificar"
$ws.Range("I56").Value = "Primera"
$ws.Range("J56").Value = 55
$ws.Range("K56").Value = 6000
$ws.Range("L56").Value = 6000
$ws.Range("M56").Value = 6000
$ws.Range("N56").Value = "$/caja 16 unidades"
$ws.Range("O56").Value = "Provincia de Quillota"
$ws.Range("P56").Value = 375
$ws.Range("Q56").Value = 16
$ws.Range("R56").Value = "Hortaliza"
